# Applies the edits described by the commit diff:
#  - On the "TestSuite" sheet: update B2, B3 ("Y" -> "N") and B4 ("N" -> "Y"),
#    and move the active selection to D13.
#  - Make "TimePageTest" the active (selected) sheet/tab of the workbook.

$wb = $excel.ActiveWorkbook

# --- TestSuite sheet: update RunMode flags -----------------------------
$wsTestSuite = $wb.Worksheets.Item("TestSuite")
$wsTestSuite.Activate()

$wsTestSuite.Range("B2").Value = "N"
$wsTestSuite.Range("B3").Value = "N"
$wsTestSuite.Range("B4").Value = "Y"

# Leave the remembered selection on this sheet at D13.
$wsTestSuite.Range("D13").Select()

# --- Make TimePageTest the active tab of the workbook -------------------
$wsTimePage = $wb.Worksheets.Item("TimePageTest")
$wsTimePage.Activate()
